# Auto-generated Excel COM-interop edit script
# Applies the BRVM Recommandations + Top_YTD automatic update

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Sheet 'Recommandations': rewrite rows 2-38 with the refreshed data, then drop the old row 39 ---

$wsReco.Cells.Item(2, 1).Value = "NEI-CEDA CI"
$wsReco.Cells.Item(2, 2).Value = 0
$wsReco.Cells.Item(2, 3).Value = 4
$wsReco.Cells.Item(2, 4).Value = 3895
$wsReco.Cells.Item(2, 5).Value = 985
$wsReco.Cells.Item(2, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(2, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(3, 1).Value = "BRVM - SERVICES PUBLICS"
$wsReco.Cells.Item(3, 2).Value = 0
$wsReco.Cells.Item(3, 3).Value = 8
$wsReco.Cells.Item(3, 4).Value = 3354.04
$wsReco.Cells.Item(3, 5).Value = 115.33
$wsReco.Cells.Item(3, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(3, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(4, 1).Value = "AIR LIQUIDE CI"
$wsReco.Cells.Item(4, 2).Value = 0
$wsReco.Cells.Item(4, 3).Value = 4
$wsReco.Cells.Item(4, 4).Value = 2795
$wsReco.Cells.Item(4, 5).Value = 700
$wsReco.Cells.Item(4, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(4, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$wsReco.Cells.Item(5, 2).Value = 0
$wsReco.Cells.Item(5, 3).Value = 4
$wsReco.Cells.Item(5, 4).Value = 2409.23
$wsReco.Cells.Item(5, 5).Value = 596.16
$wsReco.Cells.Item(5, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(5, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(6, 1).Value = "BRVM - DISTRIBUTION"
$wsReco.Cells.Item(6, 2).Value = 0
$wsReco.Cells.Item(6, 3).Value = 4
$wsReco.Cells.Item(6, 4).Value = 1993.92
$wsReco.Cells.Item(6, 5).Value = 495.1
$wsReco.Cells.Item(6, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(6, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(7, 1).Value = "BRVM - TRANSPORT"
$wsReco.Cells.Item(7, 2).Value = 0
$wsReco.Cells.Item(7, 3).Value = 4
$wsReco.Cells.Item(7, 4).Value = 1429.47
$wsReco.Cells.Item(7, 5).Value = 354.92
$wsReco.Cells.Item(7, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(7, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(8, 1).Value = "BRVM - AGRICULTURE"
$wsReco.Cells.Item(8, 2).Value = 0
$wsReco.Cells.Item(8, 3).Value = 4
$wsReco.Cells.Item(8, 4).Value = 1329.7
$wsReco.Cells.Item(8, 5).Value = 330.25
$wsReco.Cells.Item(8, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(8, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(9, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Cells.Item(9, 2).Value = 0
$wsReco.Cells.Item(9, 3).Value = 4
$wsReco.Cells.Item(9, 4).Value = 699.95
$wsReco.Cells.Item(9, 5).Value = 172.13
$wsReco.Cells.Item(9, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(9, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(10, 1).Value = "BRVM - FINANCES"
$wsReco.Cells.Item(10, 2).Value = 0
$wsReco.Cells.Item(10, 3).Value = 4
$wsReco.Cells.Item(10, 4).Value = 575.63
$wsReco.Cells.Item(10, 5).Value = 144.9
$wsReco.Cells.Item(10, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(10, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(11, 1).Value = "BRVM-PRESTIGE"
$wsReco.Cells.Item(11, 2).Value = 0
$wsReco.Cells.Item(11, 3).Value = 4
$wsReco.Cells.Item(11, 4).Value = 565.95
$wsReco.Cells.Item(11, 5).Value = 142.12
$wsReco.Cells.Item(11, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(11, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(12, 1).Value = "BRVM - SERVICES FINANCIERS"
$wsReco.Cells.Item(12, 2).Value = 0
$wsReco.Cells.Item(12, 3).Value = 4
$wsReco.Cells.Item(12, 4).Value = 565.72
$wsReco.Cells.Item(12, 5).Value = 142.41
$wsReco.Cells.Item(12, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(12, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(13, 1).Value = "BRVM - INDUSTRIELS"
$wsReco.Cells.Item(13, 2).Value = 0
$wsReco.Cells.Item(13, 3).Value = 4
$wsReco.Cells.Item(13, 4).Value = 496.97
$wsReco.Cells.Item(13, 5).Value = 122.1
$wsReco.Cells.Item(13, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(13, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(14, 1).Value = "BRVM - ENERGIE"
$wsReco.Cells.Item(14, 2).Value = 0
$wsReco.Cells.Item(14, 3).Value = 4
$wsReco.Cells.Item(14, 4).Value = 433.96
$wsReco.Cells.Item(14, 5).Value = 108.37
$wsReco.Cells.Item(14, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(14, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(15, 1).Value = "BRVM - TELECOMMUNICATIONS"
$wsReco.Cells.Item(15, 2).Value = 0
$wsReco.Cells.Item(15, 3).Value = 4
$wsReco.Cells.Item(15, 4).Value = 376.28
$wsReco.Cells.Item(15, 5).Value = 93.35
$wsReco.Cells.Item(15, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(15, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(16, 1).Value = "UNILEVER CI (UNLC)"
$wsReco.Cells.Item(16, 2).Value = 4
$wsReco.Cells.Item(16, 3).Value = 0
$wsReco.Cells.Item(16, 4).Value = 29.97
$wsReco.Cells.Item(16, 5).Value = 7.5
$wsReco.Cells.Item(16, 6).Value = "🟢 Achat"
$wsReco.Cells.Item(16, 7).Value = "✅ Renforcer"

$wsReco.Cells.Item(17, 1).Value = "SICOR CI (SICC)"
$wsReco.Cells.Item(17, 2).Value = 2
$wsReco.Cells.Item(17, 3).Value = 1
$wsReco.Cells.Item(17, 4).Value = 9.42
$wsReco.Cells.Item(17, 5).Value = 7.37
$wsReco.Cells.Item(17, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(17, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(18, 1).Value = "NESTLE CI (NTLC)"
$wsReco.Cells.Item(18, 2).Value = 2
$wsReco.Cells.Item(18, 3).Value = 0
$wsReco.Cells.Item(18, 4).Value = 8.93
$wsReco.Cells.Item(18, 5).Value = 4.58
$wsReco.Cells.Item(18, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(18, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(19, 1).Value = "SOLIBRA CI (SLBC)"
$wsReco.Cells.Item(19, 2).Value = 1
$wsReco.Cells.Item(19, 3).Value = 0
$wsReco.Cells.Item(19, 4).Value = 7.38
$wsReco.Cells.Item(19, 5).Value = 7.38
$wsReco.Cells.Item(19, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(19, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(20, 1).Value = "BERNABE CI (BNBC)"
$wsReco.Cells.Item(20, 2).Value = 1
$wsReco.Cells.Item(20, 3).Value = 0
$wsReco.Cells.Item(20, 4).Value = 5.32
$wsReco.Cells.Item(20, 5).Value = 5.32
$wsReco.Cells.Item(20, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(20, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(21, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$wsReco.Cells.Item(21, 2).Value = 2
$wsReco.Cells.Item(21, 3).Value = 1
$wsReco.Cells.Item(21, 4).Value = 4.55
$wsReco.Cells.Item(21, 5).Value = 4.35
$wsReco.Cells.Item(21, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(21, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(22, 1).Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Cells.Item(22, 2).Value = 1
$wsReco.Cells.Item(22, 3).Value = 0
$wsReco.Cells.Item(22, 4).Value = 4.45
$wsReco.Cells.Item(22, 5).Value = 4.45
$wsReco.Cells.Item(22, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(22, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(23, 1).Value = "ONATEL BF (ONTBF)"
$wsReco.Cells.Item(23, 2).Value = 1
$wsReco.Cells.Item(23, 3).Value = 0
$wsReco.Cells.Item(23, 4).Value = 3.96
$wsReco.Cells.Item(23, 5).Value = 3.96
$wsReco.Cells.Item(23, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(23, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(24, 1).Value = "BICI CI (BICC)"
$wsReco.Cells.Item(24, 2).Value = 1
$wsReco.Cells.Item(24, 3).Value = 0
$wsReco.Cells.Item(24, 4).Value = 2.05
$wsReco.Cells.Item(24, 5).Value = 2.05
$wsReco.Cells.Item(24, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(24, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(25, 1).Value = "TOTAL"
$wsReco.Cells.Item(25, 2).Value = 0
$wsReco.Cells.Item(25, 3).Value = 3
$wsReco.Cells.Item(25, 4).Value = 0
$wsReco.Cells.Item(25, 5).Value = 0
$wsReco.Cells.Item(25, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(25, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(26, 1).Value = "SETAO CI (STAC)"
$wsReco.Cells.Item(26, 2).Value = 1
$wsReco.Cells.Item(26, 3).Value = 1
$wsReco.Cells.Item(26, 4).Value = -0.72
$wsReco.Cells.Item(26, 5).Value = -4.17
$wsReco.Cells.Item(26, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(26, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(27, 1).Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Cells.Item(27, 2).Value = 1
$wsReco.Cells.Item(27, 3).Value = 2
$wsReco.Cells.Item(27, 4).Value = -2.5
$wsReco.Cells.Item(27, 5).Value = 4.29
$wsReco.Cells.Item(27, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(27, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(28, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$wsReco.Cells.Item(28, 2).Value = 0
$wsReco.Cells.Item(28, 3).Value = 1
$wsReco.Cells.Item(28, 4).Value = -2.73
$wsReco.Cells.Item(28, 5).Value = -2.73
$wsReco.Cells.Item(28, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(28, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(29, 1).Value = "SAFCA CI (SAFC)"
$wsReco.Cells.Item(29, 2).Value = 1
$wsReco.Cells.Item(29, 3).Value = 3
$wsReco.Cells.Item(29, 4).Value = -2.83
$wsReco.Cells.Item(29, 5).Value = -4.23
$wsReco.Cells.Item(29, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(29, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(30, 1).Value = "SUCRIVOIRE (SCRC)"
$wsReco.Cells.Item(30, 2).Value = 0
$wsReco.Cells.Item(30, 3).Value = 1
$wsReco.Cells.Item(30, 4).Value = -3.21
$wsReco.Cells.Item(30, 5).Value = -3.21
$wsReco.Cells.Item(30, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(30, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(31, 1).Value = "CIE CI (CIEC)"
$wsReco.Cells.Item(31, 2).Value = 0
$wsReco.Cells.Item(31, 3).Value = 1
$wsReco.Cells.Item(31, 4).Value = -3.27
$wsReco.Cells.Item(31, 5).Value = -3.27
$wsReco.Cells.Item(31, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(31, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(32, 1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$wsReco.Cells.Item(32, 2).Value = 0
$wsReco.Cells.Item(32, 3).Value = 1
$wsReco.Cells.Item(32, 4).Value = -3.51
$wsReco.Cells.Item(32, 5).Value = -3.51
$wsReco.Cells.Item(32, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(32, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(33, 1).Value = "SMB CI (SMBC)"
$wsReco.Cells.Item(33, 2).Value = 0
$wsReco.Cells.Item(33, 3).Value = 1
$wsReco.Cells.Item(33, 4).Value = -3.86
$wsReco.Cells.Item(33, 5).Value = -3.86
$wsReco.Cells.Item(33, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(33, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(34, 1).Value = "CFAO MOTORS CI (CFAC)"
$wsReco.Cells.Item(34, 2).Value = 0
$wsReco.Cells.Item(34, 3).Value = 1
$wsReco.Cells.Item(34, 4).Value = -4.97
$wsReco.Cells.Item(34, 5).Value = -4.97
$wsReco.Cells.Item(34, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(34, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(35, 1).Value = "FILTISAC CI (FTSC)"
$wsReco.Cells.Item(35, 2).Value = 0
$wsReco.Cells.Item(35, 3).Value = 2
$wsReco.Cells.Item(35, 4).Value = -6.09
$wsReco.Cells.Item(35, 5).Value = -2.89
$wsReco.Cells.Item(35, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(35, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(36, 1).Value = "NEI-CEDA CI (NEIC)"
$wsReco.Cells.Item(36, 2).Value = 0
$wsReco.Cells.Item(36, 3).Value = 1
$wsReco.Cells.Item(36, 4).Value = -7.07
$wsReco.Cells.Item(36, 5).Value = -7.07
$wsReco.Cells.Item(36, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(36, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(37, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Cells.Item(37, 2).Value = 0
$wsReco.Cells.Item(37, 3).Value = 1
$wsReco.Cells.Item(37, 4).Value = -7.3
$wsReco.Cells.Item(37, 5).Value = -7.3
$wsReco.Cells.Item(37, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(37, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(38, 1).Value = "SICABLE CI (CABC)"
$wsReco.Cells.Item(38, 2).Value = 1
$wsReco.Cells.Item(38, 3).Value = 2
$wsReco.Cells.Item(38, 4).Value = -7.47
$wsReco.Cells.Item(38, 5).Value = -7.42
$wsReco.Cells.Item(38, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(38, 7).Value = "👀 À surveiller"

# Row 39 no longer exists in the refreshed table - remove it so the sheet ends at row 38
$wsReco.Rows.Item(39).Delete()

# --- Sheet 'Top_YTD': rewrite rows 2-11 with the refreshed progression data ---

$wsYtd.Cells.Item(2, 1).Value = "BRVM - SERVICES PUBLICS"
$wsYtd.Cells.Item(2, 2).Value = 9788862.36

$wsYtd.Cells.Item(3, 1).Value = "NEI-CEDA CI"
$wsYtd.Cells.Item(3, 2).Value = 1327633.22

$wsYtd.Cells.Item(4, 1).Value = "AIR LIQUIDE CI"
$wsYtd.Cells.Item(4, 2).Value = 406940

$wsYtd.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$wsYtd.Cells.Item(5, 2).Value = 243165.15

$wsYtd.Cells.Item(6, 1).Value = "BRVM - DISTRIBUTION"
$wsYtd.Cells.Item(6, 2).Value = 128178.75

$wsYtd.Cells.Item(7, 1).Value = "BRVM - TRANSPORT"
$wsYtd.Cells.Item(7, 2).Value = 43654.63

$wsYtd.Cells.Item(8, 1).Value = "BRVM - AGRICULTURE"
$wsYtd.Cells.Item(8, 2).Value = 34864.97

$wsYtd.Cells.Item(9, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsYtd.Cells.Item(9, 2).Value = 5616.86

$wsYtd.Cells.Item(10, 1).Value = "BRVM - FINANCES"
$wsYtd.Cells.Item(10, 2).Value = 3439.05

$wsYtd.Cells.Item(11, 1).Value = "BRVM-PRESTIGE"
$wsYtd.Cells.Item(11, 2).Value = 3300.65
